$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New FSC timeline rows for three additional US sites (underhill, willowcreek, glees)
$newRows = New-Object "object[,]" 75,4
$newRows[0,0] = "underhill"; $newRows[0,1] = 44626; $newRows[0,2] = 98; $newRows[0,3] = 2022
$newRows[1,0] = "underhill"; $newRows[1,1] = 44627; $newRows[1,2] = 7; $newRows[1,3] = 2022
$newRows[2,0] = "underhill"; $newRows[2,1] = 44628; $newRows[2,2] = 15; $newRows[2,3] = 2022
$newRows[3,0] = "underhill"; $newRows[3,1] = 44629; $newRows[3,2] = 7; $newRows[3,3] = 2022
$newRows[4,0] = "underhill"; $newRows[4,1] = 44630; $newRows[4,2] = 7; $newRows[4,3] = 2022
$newRows[5,0] = "underhill"; $newRows[5,1] = 44631; $newRows[5,2] = 5; $newRows[5,3] = 2022
$newRows[6,0] = "underhill"; $newRows[6,1] = 44636; $newRows[6,2] = 100; $newRows[6,3] = 2022
$newRows[7,0] = "underhill"; $newRows[7,1] = 44637; $newRows[7,2] = 3; $newRows[7,3] = 2022
$newRows[8,0] = "underhill"; $newRows[8,1] = 44638; $newRows[8,2] = 1; $newRows[8,3] = 2022
$newRows[9,0] = "underhill"; $newRows[9,1] = 44267; $newRows[9,2] = 95; $newRows[9,3] = 2021
$newRows[10,0] = "underhill"; $newRows[10,1] = 44268; $newRows[10,2] = 90; $newRows[10,3] = 2021
$newRows[11,0] = "underhill"; $newRows[11,1] = 44270; $newRows[11,2] = 100; $newRows[11,3] = 2021
$newRows[12,0] = "underhill"; $newRows[12,1] = 44271; $newRows[12,2] = 98; $newRows[12,3] = 2021
$newRows[13,0] = "underhill"; $newRows[13,1] = 44272; $newRows[13,2] = 85; $newRows[13,3] = 2021
$newRows[14,0] = "underhill"; $newRows[14,1] = 44273; $newRows[14,2] = 80; $newRows[14,3] = 2021
$newRows[15,0] = "underhill"; $newRows[15,1] = 44274; $newRows[15,2] = 80; $newRows[15,3] = 2021
$newRows[16,0] = "underhill"; $newRows[16,1] = 44275; $newRows[16,2] = 75; $newRows[16,3] = 2021
$newRows[17,0] = "underhill"; $newRows[17,1] = 44276; $newRows[17,2] = 65; $newRows[17,3] = 2021
$newRows[18,0] = "underhill"; $newRows[18,1] = 43898; $newRows[18,2] = 97; $newRows[18,3] = 2020
$newRows[19,0] = "underhill"; $newRows[19,1] = 43899; $newRows[19,2] = 90; $newRows[19,3] = 2020
$newRows[20,0] = "underhill"; $newRows[20,1] = 43900; $newRows[20,2] = 60; $newRows[20,3] = 2020
$newRows[21,0] = "underhill"; $newRows[21,1] = 43901; $newRows[21,2] = 25; $newRows[21,3] = 2020
$newRows[22,0] = "underhill"; $newRows[22,1] = 43902; $newRows[22,2] = 25; $newRows[22,3] = 2020
$newRows[23,0] = "underhill"; $newRows[23,1] = 43903; $newRows[23,2] = 10; $newRows[23,3] = 2020
$newRows[24,0] = "underhill"; $newRows[24,1] = 43904; $newRows[24,2] = 5; $newRows[24,3] = 2020
$newRows[25,0] = "underhill"; $newRows[25,1] = 43539; $newRows[25,2] = 90; $newRows[25,3] = 2019
$newRows[26,0] = "underhill"; $newRows[26,1] = 43540; $newRows[26,2] = 90; $newRows[26,3] = 2019
$newRows[27,0] = "underhill"; $newRows[27,1] = 43542; $newRows[27,2] = 90; $newRows[27,3] = 2019
$newRows[28,0] = "underhill"; $newRows[28,1] = 43544; $newRows[28,2] = 80; $newRows[28,3] = 2019
$newRows[29,0] = "underhill"; $newRows[29,1] = 43545; $newRows[29,2] = 65; $newRows[29,3] = 2019
$newRows[30,0] = "underhill"; $newRows[30,1] = 43553; $newRows[30,2] = 95; $newRows[30,3] = 2019
$newRows[31,0] = "underhill"; $newRows[31,1] = 43555; $newRows[31,2] = 5; $newRows[31,3] = 2019
$newRows[32,0] = "underhill"; $newRows[32,1] = 43557; $newRows[32,2] = 5; $newRows[32,3] = 2019
$newRows[33,0] = "willowcreek"; $newRows[33,1] = 44663; $newRows[33,2] = 90; $newRows[33,3] = 2022
$newRows[34,0] = "willowcreek"; $newRows[34,1] = 44664; $newRows[34,2] = 75; $newRows[34,3] = 2022
$newRows[35,0] = "willowcreek"; $newRows[35,1] = 44668; $newRows[35,2] = 70; $newRows[35,3] = 2022
$newRows[36,0] = "willowcreek"; $newRows[36,1] = 44671; $newRows[36,2] = 40; $newRows[36,3] = 2022
$newRows[37,0] = "willowcreek"; $newRows[37,1] = 44672; $newRows[37,2] = 10; $newRows[37,3] = 2022
$newRows[38,0] = "willowcreek"; $newRows[38,1] = 44673; $newRows[38,2] = 2; $newRows[38,3] = 2022
$newRows[39,0] = "willowcreek"; $newRows[39,1] = 44272; $newRows[39,2] = 97; $newRows[39,3] = 2021
$newRows[40,0] = "willowcreek"; $newRows[40,1] = 44273; $newRows[40,2] = 96; $newRows[40,3] = 2021
$newRows[41,0] = "willowcreek"; $newRows[41,1] = 44274; $newRows[41,2] = 95; $newRows[41,3] = 2021
$newRows[42,0] = "willowcreek"; $newRows[42,1] = 44275; $newRows[42,2] = 92; $newRows[42,3] = 2021
$newRows[43,0] = "willowcreek"; $newRows[43,1] = 44276; $newRows[43,2] = 80; $newRows[43,3] = 2021
$newRows[44,0] = "willowcreek"; $newRows[44,1] = 44277; $newRows[44,2] = 7; $newRows[44,3] = 2021
$newRows[45,0] = "willowcreek"; $newRows[45,1] = 43923; $newRows[45,2] = 90; $newRows[45,3] = 2020
$newRows[46,0] = "willowcreek"; $newRows[46,1] = 43924; $newRows[46,2] = 70; $newRows[46,3] = 2020
$newRows[47,0] = "willowcreek"; $newRows[47,1] = 43926; $newRows[47,2] = 50; $newRows[47,3] = 2020
$newRows[48,0] = "willowcreek"; $newRows[48,1] = 43927; $newRows[48,2] = 25; $newRows[48,3] = 2020
$newRows[49,0] = "willowcreek"; $newRows[49,1] = 43928; $newRows[49,2] = 5; $newRows[49,3] = 2020
$newRows[50,0] = "willowcreek"; $newRows[50,1] = 43574; $newRows[50,2] = 97; $newRows[50,3] = 2019
$newRows[51,0] = "willowcreek"; $newRows[51,1] = 43575; $newRows[51,2] = 85; $newRows[51,3] = 2019
$newRows[52,0] = "willowcreek"; $newRows[52,1] = 43576; $newRows[52,2] = 45; $newRows[52,3] = 2019
$newRows[53,0] = "willowcreek"; $newRows[53,1] = 43577; $newRows[53,2] = 10; $newRows[53,3] = 2019
$newRows[54,0] = "glees"; $newRows[54,1] = 44715; $newRows[54,2] = 95; $newRows[54,3] = 2022
$newRows[55,0] = "glees"; $newRows[55,1] = 44716; $newRows[55,2] = 93; $newRows[55,3] = 2022
$newRows[56,0] = "glees"; $newRows[56,1] = 44717; $newRows[56,2] = 89; $newRows[56,3] = 2022
$newRows[57,0] = "glees"; $newRows[57,1] = 44718; $newRows[57,2] = 83; $newRows[57,3] = 2022
$newRows[58,0] = "glees"; $newRows[58,1] = 44719; $newRows[58,2] = 75; $newRows[58,3] = 2022
$newRows[59,0] = "glees"; $newRows[59,1] = 44720; $newRows[59,2] = 70; $newRows[59,3] = 2022
$newRows[60,0] = "glees"; $newRows[60,1] = 44721; $newRows[60,2] = 60; $newRows[60,3] = 2022
$newRows[61,0] = "glees"; $newRows[61,1] = 44722; $newRows[61,2] = 45; $newRows[61,3] = 2022
$newRows[62,0] = "glees"; $newRows[62,1] = 44723; $newRows[62,2] = 25; $newRows[62,3] = 2022
$newRows[63,0] = "glees"; $newRows[63,1] = 44724; $newRows[63,2] = 15; $newRows[63,3] = 2022
$newRows[64,0] = "glees"; $newRows[64,1] = 44725; $newRows[64,2] = 10; $newRows[64,3] = 2022
$newRows[65,0] = "glees"; $newRows[65,1] = 44726; $newRows[65,2] = 7; $newRows[65,3] = 2022
$newRows[66,0] = "glees"; $newRows[66,1] = 44727; $newRows[66,2] = 5; $newRows[66,3] = 2022
$newRows[67,0] = "glees"; $newRows[67,1] = 44728; $newRows[67,2] = 3; $newRows[67,3] = 2022
$newRows[68,0] = "glees"; $newRows[68,1] = 44350; $newRows[68,2] = 93; $newRows[68,3] = 2021
$newRows[69,0] = "glees"; $newRows[69,1] = 44351; $newRows[69,2] = 85; $newRows[69,3] = 2021
$newRows[70,0] = "glees"; $newRows[70,1] = 44352; $newRows[70,2] = 80; $newRows[70,3] = 2021
$newRows[71,0] = "glees"; $newRows[71,1] = 44353; $newRows[71,2] = 70; $newRows[71,3] = 2021
$newRows[72,0] = "glees"; $newRows[72,1] = 44354; $newRows[72,2] = 60; $newRows[72,3] = 2021
$newRows[73,0] = "glees"; $newRows[73,1] = 44355; $newRows[73,2] = 45; $newRows[73,3] = 2021
$newRows[74,0] = "glees"; $newRows[74,1] = 44356; $newRows[74,2] = 25; $newRows[74,3] = 2021

$ws.Range("A1204:D1278").Value = $newRows

# Update selection to match the final cursor position after the paste
[void]$ws.Range("J1278").Select()
